$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the new test run results
$ws.Range("A2").Value = "dmoralesr"
$ws.Range("C2").Value = "'2103764"
$ws.Range("F2").Value = "AAACT231729MLH238M "
$ws.Range("G2").Value = "21 jun. 2023, 14:50:44"
$ws.Range("H2").Value = 1010824482

# Update the active selection to H2 as recorded in the saved view
$ws.Range("H2").Select()
